$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 11:03"

# --- Update country stats rows whose totals changed but country/rank stayed put ---
$ws.Range("B31").Value = 16268
$ws.Range("C31").Value = 22
$ws.Range("D31").Value = 10223
$ws.Range("E31").Value = 5808
$ws.Range("F31").Value = 89
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 237

$ws.Range("B33").Value = 15650
$ws.Range("C33").Value = 29
$ws.Range("D33").Value = 13462
$ws.Range("E33").Value = 1582
$ws.Range("F33").Value = 104
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 606

$ws.Range("B36").Value = 14242
$ws.Range("C36").Value = 236
$ws.Range("E36").Value = 9262
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 700

# --- Rows 40-44: countries reshuffled ranking (Banglades overtakes Corea del Sur
#     and Filipinas climbs above Dinamarca/Serbia) together with new totals ---
$ws.Range("A40").Value = "Banglades"
$ws.Range("B40").Value = 10929
$ws.Range("C40").Value = 786
$ws.Range("D40").Value = 1403
$ws.Range("E40").Value = 9343
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 183

$ws.Range("A41").Value = "Corea del Sur"
$ws.Range("B41").Value = 10804
$ws.Range("C41").Value = 3
$ws.Range("D41").Value = 9283
$ws.Range("E41").Value = 1267
$ws.Range("F41").Value = 55
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 254

$ws.Range("A42").Value = "Filipinas"
$ws.Range("B42").Value = 9684
$ws.Range("C42").Value = 199
$ws.Range("D42").Value = 1408
$ws.Range("E42").Value = 7639
$ws.Range("F42").Value = 31
$ws.Range("G42").Value = 14
$ws.Range("H42").Value = 637

$ws.Range("A43").Value = "Dinamarca"
$ws.Range("B43").Value = 9670
$ws.Range("D43").Value = 7088
$ws.Range("E43").Value = 2089
$ws.Range("F43").Value = 57
$ws.Range("H43").Value = 493

$ws.Range("A44").Value = "Serbia"
$ws.Range("B44").Value = 9557
$ws.Range("D44").Value = 1574
$ws.Range("E44").Value = 7786
$ws.Range("F44").Value = 53
$ws.Range("H44").Value = 197

$ws.Range("B53").Value = 6383
$ws.Range("C53").Value = 30
$ws.Range("D53").Value = 4567
$ws.Range("E53").Value = 1710
$ws.Range("F53").Value = 24
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 106

# --- Rows 67-68: Oman overtakes Ghana ---
$ws.Range("A67").Value = "Oman"
$ws.Range("B67").Value = 2735
$ws.Range("C67").Value = 98
$ws.Range("D67").Value = 858
$ws.Range("E67").Value = 1865
$ws.Range("F67").Value = 17
$ws.Range("H67").Value = 12

$ws.Range("A68").Value = "Ghana"
$ws.Range("B68").Value = 2719
$ws.Range("D68").Value = 294
$ws.Range("E68").Value = 2407
$ws.Range("F68").Value = 4
$ws.Range("H68").Value = 18

$ws.Range("D72").Value = 1440
$ws.Range("E72").Value = 739

$ws.Range("F88").Value = 6

$ws.Range("D99").Value = 197
$ws.Range("E99").Value = 549
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 9

$ws.Range("B115").Value = 557
$ws.Range("C115").Value = 16
$ws.Range("E115").Value = 539
